$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: insert new HT correct-score odds columns ---
# AX1:BC1 (Odd_CS_0-1_HT .. Odd_CS_2-3_HT) shift left into AW1:BB1;
# the old AW1 value (Odd_CS_3-3_HT) moves to BC1; a new BD1 column
# (Odd_CS_4-4_HT) is appended, matching the header style (s="1").
$ws.Range("AW1").Value = "Odd_CS_0-1_HT"
$ws.Range("AX1").Value = "Odd_CS_0-2_HT"
$ws.Range("AY1").Value = "Odd_CS_1-2_HT"
$ws.Range("AZ1").Value = "Odd_CS_0-3_HT"
$ws.Range("BA1").Value = "Odd_CS_1-3_HT"
$ws.Range("BB1").Value = "Odd_CS_2-3_HT"
$ws.Range("BC1").Value = "Odd_CS_3-3_HT"
$ws.Range("BD1").Value = "Odd_CS_4-4_HT"

# New BD1 header cell needs the same bold/border/centered style as the
# rest of row 1; copy formatting from an existing header cell (A1).
$ws.Range("A1").Copy()
$ws.Range("BD1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: replace the match record with the new one, plus the new BD2 odd ---
$ws.Range("A2").Value = "GbmHWyQ7"
$ws.Range("B2").Value = "31/10/2024"
$ws.Range("C2").Value = "11:00"
$ws.Range("D2").Value = "EGYPT - PREMIER LEAGUE"
$ws.Range("E2").Value = "El Gouna"
$ws.Range("F2").Value = "ZED"
$ws.Range("G2").Value = 4.15
$ws.Range("H2").Value = 2.85
$ws.Range("I2").Value = 2.05
$ws.Range("J2").Value = 4.55
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 2.57
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 5.7
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.6
$ws.Range("Q2").Value = 2.32
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.6
$ws.Range("U2").Value = 1.93
$ws.Range("V2").Value = 1.78
$ws.Range("W2").Value = 9.5
$ws.Range("X2").Value = 22
$ws.Range("Y2").Value = 13.5
$ws.Range("Z2").Value = 75
$ws.Range("AA2").Value = 45
$ws.Range("AB2").Value = 55
$ws.Range("AC2").Value = 5.7
$ws.Range("AD2").Value = 5.6
$ws.Range("AE2").Value = 15
$ws.Range("AF2").Value = 90
$ws.Range("AG2").Value = 800
$ws.Range("AH2").Value = 6
$ws.Range("AI2").Value = 9
$ws.Range("AJ2").Value = 8.5
$ws.Range("AK2").Value = 19.5
$ws.Range("AL2").Value = 18.5
$ws.Range("AM2").Value = 32
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 24
$ws.Range("AP2").Value = 28
$ws.Range("AQ2").Value = 150
$ws.Range("AR2").Value = 175
$ws.Range("AS2").Value = 350
$ws.Range("AT2").Value = 2.6
$ws.Range("AU2").Value = 6.7
$ws.Range("AV2").Value = 60
$ws.Range("AW2").Value = 3.85
$ws.Range("AX2").Value = 10.25
$ws.Range("AY2").Value = 18
$ws.Range("AZ2").Value = 40
$ws.Range("BA2").Value = 70
$ws.Range("BB2").Value = 250
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51
